$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column D (Parts column moves from D to E, etc.)
$ws.Range("D1").EntireColumn.Insert()

# Set the new column D's width to match column C's width (~12.86 chars)
$ws.Range("D1").ColumnWidth = 12

# Fill in the new "v" / "v 4.7k" values for rows 2-11 (the voltage rating column)
$ws.Range("D2").Value = "v"
$ws.Range("D3").Value = "v"
$ws.Range("D4").Value = "v"
$ws.Range("D5").Value = "v"
$ws.Range("D6").Value = "v 4.7k"
$ws.Range("D7").Value = "v"
$ws.Range("D8").Value = "v"
$ws.Range("D9").Value = "v"
$ws.Range("D10").Value = "v"
$ws.Range("D11").Value = "v"

# Move the selection to where it ended up in the authored file
$ws.Range("D12").Select() | Out-Null
